$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update subtitle (row 2) to reflect new reporting month
$ws.Range("A2").Value = "Total (All Sectors) by Census Division and State, Year-to-Date through November 2016"

# Update data values (RSE percentages) per EPM_2016_11 run
$ws.Range("B4").Value = 4
$ws.Range("H4").Value = 14
$ws.Range("C5").Value = 90
$ws.Range("E5").Value = 3
$ws.Range("H5").Value = 93
$ws.Range("C6").Value = 108
$ws.Range("E6").Value = 9
$ws.Range("H6").Value = 18
$ws.Range("B7").Value = 4
$ws.Range("E7").Value = 5
$ws.Range("H7").Value = 39
$ws.Range("C8").Value = 61
$ws.Range("E8").Value = 1
$ws.Range("H8").Value = 26
$ws.Range("C9").Value = 28
$ws.Range("E9").Value = 1
$ws.Range("H9").Value = 989
$ws.Range("C10").Value = 537
$ws.Range("E10").Value = 258
$ws.Range("H10").Value = 45
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 33
$ws.Range("D11").Value = 69
$ws.Range("H11").Value = 2
$ws.Range("C12").Value = 206
$ws.Range("D12").Value = 128
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 75
$ws.Range("H12").Value = 886
$ws.Range("C13").Value = 67
$ws.Range("E13").Value = 3
$ws.Range("H13").Value = 2
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = 25
$ws.Range("D14").Value = 81
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 28
$ws.Range("H14").Value = 17
$ws.Range("B15").Value = 0.46
$ws.Range("C15").Value = 4
$ws.Range("F15").Value = 12
$ws.Range("H15").Value = 15
$ws.Range("B16").Value = 0.35
$ws.Range("C16").Value = 6
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 99
$ws.Range("H16").Value = 93
$ws.Range("C17").Value = 6
$ws.Range("F17").Value = 18
$ws.Range("H17").Value = 15
$ws.Range("C18").Value = 15
$ws.Range("D18").Value = 22
$ws.Range("H18").Value = 30
$ws.Range("D19").Value = 4
$ws.Range("F19").Value = 61
$ws.Range("H19").Value = 26
$ws.Range("B20").Value = 0.21
$ws.Range("C20").Value = 29
$ws.Range("E20").Value = 2
$ws.Range("H20").Value = 24
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 152
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 100
$ws.Range("H21").Value = 8
$ws.Range("C22").Value = 19
$ws.Range("D22").Value = 152
$ws.Range("E22").Value = 15
$ws.Range("H22").Value = 37
$ws.Range("C23").Value = 26
$ws.Range("E23").Value = 14
$ws.Range("H23").Value = 354
$ws.Range("C24").Value = 140
$ws.Range("H24").Value = 38
$ws.Range("E25").Value = 12
$ws.Range("H25").Value = 13
$ws.Range("C26").Value = 128
$ws.Range("E26").Value = 8
$ws.Range("H26").Value = 28
$ws.Range("C27").Value = 12
$ws.Range("E27").Value = 42
$ws.Range("F27").Value = 100
$ws.Range("C28").Value = 508
$ws.Range("E28").Value = 14
$ws.Range("H28").Value = 0.44
$ws.Range("B29").Value = 0.16
$ws.Range("C29").Value = 9
$ws.Range("E29").Value = 0.48
$ws.Range("C30").Value = 367
$ws.Range("E30").Value = 7
$ws.Range("C31").Value = 3091
$ws.Range("E31").Value = 156
$ws.Range("B32").Value = 0.32
$ws.Range("C32").Value = 8
$ws.Range("E32").Value = 1
$ws.Range("H32").Value = 92
$ws.Range("C33").Value = 53
$ws.Range("H33").Value = 15
$ws.Range("C34").Value = 38
$ws.Range("E34").Value = 24
$ws.Range("H34").Value = 0
$ws.Range("C35").Value = 33
$ws.Range("H35").Value = 17
$ws.Range("B36").Value = 0.04
$ws.Range("C36").Value = 57
$ws.Range("E36").Value = 1
$ws.Range("H36").Value = 27
$ws.Range("C37").Value = 36
$ws.Range("H37").Value = 38
$ws.Range("B38").Value = 0.2
$ws.Range("E38").Value = 3
$ws.Range("H38").Value = 30
$ws.Range("B39").Value = 0.27
$ws.Range("C39").Value = 6
$ws.Range("F39").Value = 72
$ws.Range("H39").Value = 7
$ws.Range("B40").Value = 0.28999999999999998
$ws.Range("C40").Value = 45
$ws.Range("E40").Value = 2
$ws.Range("F40").Value = 119
$ws.Range("H40").Value = 11
$ws.Range("B41").Value = 1
$ws.Range("E41").Value = 3
$ws.Range("H41").Value = 8
$ws.Range("C42").Value = 388
$ws.Range("C43").Value = 2
$ws.Range("E43").Value = 2
$ws.Range("H43").Value = 12
$ws.Range("B44").Value = 0.08
$ws.Range("C44").Value = 5
$ws.Range("D44").Value = 3
$ws.Range("E44").Value = 1
$ws.Range("H44").Value = 12
$ws.Range("C45").Value = 3
$ws.Range("E45").Value = 1
$ws.Range("H45").Value = 17
$ws.Range("D46").Value = 2
$ws.Range("F46").Value = 7
$ws.Range("B47").Value = 1
$ws.Range("C47").Value = 12
$ws.Range("H47").Value = 30
$ws.Range("C48").Value = 14
$ws.Range("D48").Value = 59
$ws.Range("F48").Value = 5
$ws.Range("H48").Value = 41
$ws.Range("C49").Value = 12
$ws.Range("F49").Value = 9
$ws.Range("H49").Value = 4
$ws.Range("C50").Value = 9
$ws.Range("E50").Value = 0.4
$ws.Range("B51").Value = 0.18
$ws.Range("C51").Value = 87
$ws.Range("H51").Value = 33
$ws.Range("B52").Value = 105
$ws.Range("C52").Value = 475
$ws.Range("E52").Value = 9
$ws.Range("H52").Value = 10
$ws.Range("C53").Value = 51
$ws.Range("E53").Value = 59
$ws.Range("H53").Value = 4
$ws.Range("C55").Value = 56
$ws.Range("E55").Value = 5
$ws.Range("H55").Value = 116
$ws.Range("C56").Value = 28
$ws.Range("E56").Value = 7
$ws.Range("F56").Value = 457
$ws.Range("H56").Value = 41
$ws.Range("C57").Value = 2
$ws.Range("E57").Value = 22
$ws.Range("F57").Value = 7
$ws.Range("H57").Value = 28
$ws.Range("C58").Value = 37
$ws.Range("F58").Value = 6
$ws.Range("C59").Value = 32
$ws.Range("E59").Value = 2
$ws.Range("F59").Value = 8
$ws.Range("H59").Value = 7
$ws.Range("C60").Value = 13658
$ws.Range("C61").Value = 80
$ws.Range("B62").Value = 5
$ws.Range("C62").Value = 6
$ws.Range("E62").Value = 13
$ws.Range("F62").Value = 130
$ws.Range("H62").Value = 20
$ws.Range("B63").Value = 18
$ws.Range("C63").Value = 6
$ws.Range("E63").Value = 13
$ws.Range("H63").Value = 20
$ws.Range("B64").Value = 3
$ws.Range("C64").Value = 6
$ws.Range("F64").Value = 130
$ws.Range("H64").Value = 79
$ws.Range("C65").Value = 4
$ws.Range("D65").Value = 3
$ws.Range("E65").Value = 0.34
$ws.Range("F65").Value = 6
